# Commit: #documentSource Update document of source
#
# Appends eight new paragraphs to the end of the document, right after
# the "SystemDB: Dữ liệu hệ thống" bullet, describing the Services /
# ViewModels source layout. Built as a WordprocessingML fragment and
# dropped in via Range.InsertXML so the new paragraphs get exactly the
# pPr/rPr shape (numbered ListParagraph bullets + plain indented
# "section title" paragraphs, ending in one blank indented paragraph)
# that the target document has.

$d = $word.ActiveDocument

# Collapsed range sitting at the very end of the document body (right
# after the last run of the last existing paragraph).
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$insertXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Những phần ở ngoài thường sẽ là phần quan trọng, dường như là một module cho hệ thống nếu cần bổ sung gì khác lớn</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="567" w:firstLine="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Services</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Là những hàm chức năng không thuộc model nhưng cần phải có như genare token, check token.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="567" w:firstLine="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ViewModels</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Những chức năng của model sẽ ở đây thay vì như ở service như các phần mềm khác</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Dựa vào Models mà tạo tổ chức các tập tin và thư mục khác nhau</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Những chức năng trong đây sẽ liên quan mật thiết đến model như thêm mới, cập nhật, xóa, kiểm tra tồn tại dữ liệu, đếm dữ liệu…</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="567" w:firstLine="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($insertXml)
